$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are managing a program that requires co-ordinating with multiple stakeholders, including external vendors. However, an external stakeholder is missing their deadlines, creating overall delays.What action should you take?",
        "ques_type": 2,
        "options": [
            "Inform the stakeholder that they need to get on track immediately.",
            "Inform the stakeholder of the impact their delay is causing.",
            "Create a mitigation plan in discussion with the external stakeholder.",
            "Escalate the issue to the stakeholder\u2019s manager and ask for their assistance."
        ],
        "score": "Inform the stakeholder of the impact their delay is causing."
    },
    {
        "title": "As a program manager leading a team of software engineers, you have recently assigned tasks. A highly competent engineer is frustrated that they have been assigned a task they describe as \u201crepetitive and boring\u201d and has requested you assign them a different task. However, they are the only member of the team with the knowledge required to complete it. What action should you take?",
        "ques_type": 2,
        "options": [
            "Empathize with the engineer, but ask them to complete the task explaining the critical nature of the job.",
            "Work with the engineer to identify whether the creation of a transition plan would be feasible.",
            "Refuse the engineer\u2019s request, but promise them an additional bonus for continuing to work on the task.",
            "Commit to considering the engineer\u2019s request at the end of the quarter, explaining the immediate importance of the task."
        ],
        "score": "Empathize with the engineer, but ask them to complete the task explaining the critical nature of the job."
    },
    {
        "title": "Your product manager has asked for a new feature to be built. However, you and your team feel this feature would not benefit the users, and efforts should be focused on other features in the backlog that will have a bigger impact. What action should you take?",
        "ques_type": 2,
        "options": [
            "Implement the new feature.",
            "Refuse to work on the new feature.",
            "Explain your rationale to your product manager in an attempt to reach a consensus.",
            "Request your product manager to reconsider because there is consensus within the team on priority."
        ],
        "score": "Implement the new feature."
    },
    {
        "title": "You and a colleague are attending a leadership off-site to discuss your product\u2019s roadmap. In discussions with other leaders, you have both started to realize that few of them are aware of the program you are leading. Your colleague has just asked you what you should do about it.What should you tell them?",
        "ques_type": 2,
        "options": [
            "\u201cGive it time. The great work we are doing will speak for itself.\u201d",
            "\u201cStart sharing details of our program with everyone you speak to here.\u201d",
            "\u201cLet\u2019s create a communication plan for how to keep stakeholders updated.\u201d",
            "\u201cI\u2019m going to set up regular one-on-one meetings with my manager so I can share updates.\u201d"
        ],
        "score": "\u201cStart sharing details of our program with everyone you speak to here.\u201d"
    }
]
'@

# The old layout had A1 = 0 (bold, bordered, centered style) and A2 =
# the question-bank text (default style). Put the reformatted text in
# A2, then cut/paste it into A1 (preserving its default formatting)
# and remove the now-empty old row 2.
$ws.Range("A2").Value = $text
$ws.Range("A2").Cut($ws.Range("A1"))
$ws.Rows.Item(2).Delete()
